$d = $word.ActiveDocument

# --- Step 1: remove the bookmarkStart/bookmarkEnd (_GoBack) from paragraph 1 ---
$p1 = $d.Paragraphs.Item(1)
$x1 = $p1.Range.WordOpenXML
if ($x1 -notmatch '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>') {
    throw "paragraph1: _GoBack bookmark pair not found"
}
$x1 = $x1 -replace '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>', ''
$p1.Range.InsertXML($x1)

# --- Step 2: paragraph with the two "Julio"/date WordArt pictures ---
$p5 = $d.Paragraphs.Item($d.Paragraphs.Count)
$x5 = $p5.Range.WordOpenXML

# 2a. first textpath: "Julio        2023<LF>" -> "Agosto        2023<LF>", drop the w10:border* quad
$old1 = '<v:textpath style="font-family:&quot;Impact&quot;;font-size:18pt;v-text-kern:t" trim="t" fitpath="t" xscale="f" string="Julio        2023&#xA;"/><w10:bordertop type="single" width="4"/><w10:borderleft type="single" width="4"/><w10:borderbottom type="single" width="4"/><w10:borderright type="single" width="4"/>'
$new1 = '<v:textpath style="font-family:&quot;Impact&quot;;font-size:18pt;v-text-kern:t" trim="t" fitpath="t" xscale="f" string="Agosto        2023&#xA;"/>'
if ($x5.IndexOf($old1) -lt 0) { throw "paragraph5: first textpath pattern not found" }
$x5 = $x5.Replace($old1, $new1)

# 2b. second textpath: "Julio   2023" -> "Agosto   2023", drop the w10:border* quad
$old2 = '<v:textpath style="font-family:&quot;Impact&quot;;font-size:18pt;v-text-kern:t" trim="t" fitpath="t" xscale="f" string="Julio   2023"/><w10:bordertop type="single" width="4"/><w10:borderleft type="single" width="4"/><w10:borderbottom type="single" width="4"/><w10:borderright type="single" width="4"/>'
$new2 = '<v:textpath style="font-family:&quot;Impact&quot;;font-size:18pt;v-text-kern:t" trim="t" fitpath="t" xscale="f" string="Agosto   2023"/>'
if ($x5.IndexOf($old2) -lt 0) { throw "paragraph5: second textpath pattern not found" }
$x5 = $x5.Replace($old2, $new2)

# 2c. move the _GoBack bookmark: starts right before the yellow-highlight run holding
#     the second picture, ends right after that same run (before </w:p>)
$oldMarker = '<w:r><w:tab/></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:pict w14:anchorId="0AF744A5">'
if ($x5.IndexOf($oldMarker) -lt 0) { throw "paragraph5: bookmark-insertion marker not found" }
$newMarker = '<w:r><w:tab/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:pict w14:anchorId="0AF744A5">'
$x5 = $x5.Replace($oldMarker, $newMarker)

$oldEnd = '</w:pict></w:r></w:p>'
if ($x5.IndexOf($oldEnd) -lt 0) { throw "paragraph5: end-of-paragraph marker not found" }
$newEnd = '</w:pict></w:r><w:bookmarkEnd w:id="0"/></w:p>'
$x5 = $x5.Replace($oldEnd, $newEnd)

$p5.Range.InsertXML($x5)

Write-Output "done"
